# Commit: "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab Neo4j query stored in cell B4 of the "startup" sheet is
# corrected: the `File Type` and `Breed` columns are dropped from the
# RETURN clause (those properties aren't reliably present for every file,
# so the script was edited to stop requesting them).
#
# Because the text got shorter, Excel's wrapped-row auto-height for row 4
# shrinks from 246.5 to 217.5. The active selection also moves from C2
# (the CasesTab script cell) onto the freshly edited B4 cell, scrolling
# the view so row 4 is visible at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['English Setter']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# B4 holds the FilesTab script (s="1" -> wrap-text style already applied).
$ws.Range("B4").Value = $newFilesQuery

# Shorter text -> smaller auto-fit row height for the wrapped cell.
$ws.Rows.Item(4).RowHeight = 217.5

# Move the selection/view onto the cell that was just edited.
[void]$ws.Range("A4").Select()
[void]$ws.Range("B4").Select()
